$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading 1) paragraph at the top of the document:
#       [empty run]
#       [bold]     "Meta description"
#       [normal]   ": Read our 2021 review of 3Diamonds FashionTV and
#                    play for free. Discover the special features and
#                    bonuses, pros and cons, and RTP and volatility
#                    values."
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">: Read our 2021 review of 3Diamonds FashionTV and play for free. Discover the special features and bonuses, pros and cons, and RTP and volatility values.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $metaPara.Range
$insertPoint.Collapse(1)
$insertPoint.InsertXML($metaXml)

Write-Output "Inserted meta description paragraph; new paragraph count = $($d.Paragraphs.Count)"

# ------------------------------------------------------------------
# 2) Near the end of the document there are (still) two paragraphs:
#       [bold]   "Play 3Diamonds FashionTV Slot for Free - Review 2021"
#       [italic] "Read our 2021 review of 3Diamonds FashionTV and play
#                 for free. Discover the special features and bonuses,
#                 pros and cons, and RTP and volatility values."
#    The bold title paragraph (a duplicate of the doc's first heading)
#    is removed entirely, and the italic paragraph's text is rewritten
#    into an image-generation prompt, keeping its italic formatting.
# ------------------------------------------------------------------
$oldTitleText = "Play 3Diamonds FashionTV Slot for Free - Review 2021"
$oldMetaText = "Read our 2021 review of 3Diamonds FashionTV and play for free. Discover the special features and bonuses, pros and cons, and RTP and volatility values."
$newImagePrompt = 'Please create a cartoon style feature image for "3Diamonds FashionTV" that features a Maya warrior wearing glasses and looking happy. The image should be fun and vibrant, incorporating elements of luxury such as champagne bottles, diamonds, credit cards, and cars. Additionally, the image should highlight the nighttime city backdrop with a golden glow. Try to capture the excitement and uniqueness of this slot game with your image and make it stand out to potential players.'

$count = $d.Paragraphs.Count
$bottomTitlePara = $d.Paragraphs($count - 1)
$bottomMetaPara = $d.Paragraphs($count)

if ($bottomTitlePara.Range.Text.TrimEnd() -eq $oldTitleText -and $bottomMetaPara.Range.Text.TrimEnd() -eq $oldMetaText) {
    # Remove the duplicated bold title paragraph entirely (its paragraph
    # mark goes with it, merging cleanly with what follows).
    $bottomTitlePara.Range.Delete()

    $count = $d.Paragraphs.Count
    $targetPara = $d.Paragraphs($count)
    $targetRange = $targetPara.Range
    # Exclude the trailing paragraph mark so only the run text changes
    # (this keeps the existing italic run/formatting intact and avoids
    # Find/Replace's smart-quote autocorrection).
    $targetRange.MoveEnd(1, -1)
    $targetRange.Text = $newImagePrompt

    Write-Output "Replaced trailing paragraphs; final paragraph count = $($d.Paragraphs.Count)"
} else {
    Write-Output "WARNING: trailing paragraphs did not match expected text; no changes made there."
    Write-Output "bottomTitlePara=[$($bottomTitlePara.Range.Text)]"
    Write-Output "bottomMetaPara=[$($bottomMetaPara.Range.Text)]"
}
